$d = $word.ActiveDocument

# ------------------------------------------------------------------------
# 1) The paragraph "DESCRIÇÃO: CNPJ, Nome Fantasia, Endereço, Estado,
#    Telefone." gets its last sentence edited: ", Cep" is typed in just
#    before the final period. Word records the edit location with its
#    auto-managed "_GoBack" bookmark, which now wraps the freshly
#    (re)typed text, and the sentence ends up split into three runs
#    ("...Telefone" / ", Cep" / ".") because of how the text was typed.
# ------------------------------------------------------------------------
$cnpjFind = $d.Content.Duplicate
$cnpjFind.Find.Execute("DESCRIÇÃO: CNPJ, Nome Fantasia, Endereço, Estado, Telefone.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$cnpjMatch = $d.Range($cnpjFind.Start, $cnpjFind.End)
$cnpjPara = $cnpjMatch.Paragraphs(1).Range
$cnpjParagraph = $d.Range($cnpjPara.Start, $cnpjPara.End)

$cnpjXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p w:rsidR="00C57325" w:rsidRPr="00C57325" w:rsidRDefault="00C57325" w:rsidP="00020DC5">
            <w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/></w:rPr></w:pPr>
            <w:r w:rsidRPr="00C57325"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/></w:rPr><w:t>DESCRIÇÃO:</w:t></w:r>
            <w:r w:rsidRPr="00C57325"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>CNPJ, Nome Fantasia, Endereço, Estado, Telefone</w:t></w:r>
            <w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>, Cep</w:t></w:r>
            <w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>.</w:t></w:r>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$cnpjParagraph.InsertXML($cnpjXml)

# ------------------------------------------------------------------------
# 2) Because "_GoBack" can only mark one spot in the document, the pair
#    of (now empty) bookmarkStart/bookmarkEnd tags that used to sit right
#    after "RINF10" - marking the *previous* last-edit location - is
#    removed.
# ------------------------------------------------------------------------
$rinfFind = $d.Content.Duplicate
$rinfFind.Find.Execute("ID: RINF10", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rinfMatch = $d.Range($rinfFind.Start, $rinfFind.End)
$rinfPara = $rinfMatch.Paragraphs(1).Range
$rinfParagraph = $d.Range($rinfPara.Start, $rinfPara.End)

$rinfXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p w:rsidR="00010367" w:rsidRPr="00010367" w:rsidRDefault="00010367" w:rsidP="00755DD7">
            <w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr>
            <w:r w:rsidRPr="00010367"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/></w:rPr><w:t>ID:</w:t></w:r>
            <w:r w:rsidRPr="00010367"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
            <w:r w:rsidR="00D7782E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>RINF10</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$rinfParagraph.InsertXML($rinfXml)
